# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
# per the diff: several rows' F values increased.

$wb = $excel.ActiveWorkbook

# Map of row -> new F value for sheet "展览" (sheet1)
$updatesExhibition = @{
    2  = 175
    5  = 1285
    6  = 17944
    7  = 352
    8  = 256
    9  = 1065
    10 = 6769
    11 = 681
    15 = 62
    17 = 147
    18 = 1299
    19 = 204
    20 = 55
    24 = 33
    26 = 974
    27 = 110
    28 = 5158
    30 = 12
    31 = 8
    33 = 11986
    34 = 1272
    35 = 39
    36 = 200
    39 = 297
}

# Map of row -> new F value for sheet "全部类型" (sheet4)
$updatesAllTypes = @{
    2  = 175
    5  = 1285
    6  = 17944
    7  = 352
    8  = 256
    9  = 1065
    10 = 6769
    11 = 681
    15 = 62
    17 = 147
    18 = 1299
    19 = 204
    20 = 55
    24 = 33
    26 = 974
    27 = 110
    28 = 5158
    32 = 12
    33 = 8
    35 = 11986
    36 = 1272
    37 = 39
    38 = 200
    41 = 297
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $wsExhibition.Range("F$row").Value = $updatesExhibition[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAllTypes.Keys) {
    $wsAllTypes.Range("F$row").Value = $updatesAllTypes[$row]
}
